$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds the (mis-dated) game date as text, e.g. "5-23-2007-08".
# Rewrite it to the correct ISO-ish date string "2008-05-23" for every
# data row (rows 2-31; row 1 is the "Date" header).
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    # A leading apostrophe forces Excel to keep this as literal text
    # instead of auto-converting the ISO-looking string into a date
    # serial number.
    $cell.Value = "'2008-05-23"
    # Re-apply the (unstyled) Normal cell style so this text entry
    # doesn't pick up a distinct "quote prefix" style from the literal
    # text assignment above - the cell had no explicit style before.
    $cell.Style = "Normal"
}
